$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New nitrate mass-balance block (columns P:T, rows 18-30) ---
# String-valued cells are written first, in the exact order the source
# workbook introduced them, so new shared-string entries land at the same
# indices as the authored file.
$ws.Range("P19").Value = "initial ext"
$ws.Range("P20").Value = "initial int"
$ws.Range("Q20").Value = "x"
$ws.Range("R19").Value = "mmol/l"
$ws.Range("R20").Value = "mmol/gDW"
$ws.Range("P22").Value = "final ext"
$ws.Range("P24").Value = "final_biomass"
$ws.Range("P25").Value = "comp"
$ws.Range("R24").Value = " g/L"
$ws.Range("P27").Value = "nitrate in biomass"
$ws.Range("R27").Value = "mmol/L"
$ws.Range("T20").Value = "y"

# Citation link added on row 5 of the first Km/Vmax block (written last,
# matching shared-string index 20 in the authored file).
$ws.Range("E5").Value = "https://www.sciencedirect.com/science/article/pii/S0960852414012486#ab005"

# Remaining (non-string) cells of the nitrate block.
$ws.Range("Q19").Value = 10
$ws.Range("T19").Value = 7.5

$ws.Range("R22").Value = "mmol/l"
$ws.Range("Q22").Formula = "=10-7.5"
$ws.Range("T22").Value = 0

$ws.Range("Q24").Value = 1.024
$ws.Range("T24").Value = 1.76

$ws.Range("R25").Value = "mmol/gDW"
$ws.Range("Q25").Value = 4.42
$ws.Range("T25").Value = 4.42

$ws.Range("Q27").Formula = "=Q24*Q25"
$ws.Range("T27").Formula = "=T25*T24"

$ws.Range("T29").Formula = "=T27-T19"
$ws.Range("T30").Formula = "=T29/T18"

$ws.Range("T18").Value = 0.11799999999999999

# --- New weighted-average helper cell ---
$ws.Range("S6").Formula = "=11.95*0.25+18*0.1+11.05*0.45+8*0.2"

# --- D11 average now also includes the E11:F11 pair ---
$ws.Range("D11").Formula = "=AVERAGE(B11:C11,E11:F11)"

# --- Selection to mirror the saved workbook ---
[void]$ws.Range("E6").Select()
